$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 111870313
$ws.Range("B2").Value = 90800
$ws.Range("D2").Value = "LC"
$ws.Range("E2").Value = 4364
$ws.Range("F2").Value = "Dropptaggsvamp"
$ws.Range("G2").Value = "Hydnellum ferrugineum"
$ws.Range("H2").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q2").Value = 494301
$ws.Range("R2").Value = 6928922
$ws.Range("S2").Value = 20
# Row 3
$ws.Range("A3").Value = 111868481
$ws.Range("B3").Value = 89902
$ws.Range("D3").Value = "VU"
$ws.Range("E3").Value = 298
$ws.Range("F3").Value = "Laxgröppa"
$ws.Range("G3").Value = "Byssomerulius albostramineus"
$ws.Range("H3").Value = "(Torrend) Hjortstam"
$ws.Range("Q3").Value = 494354
$ws.Range("R3").Value = 6928891
# Row 4
$ws.Range("A4").Value = 111868443
$ws.Range("B4").Value = 94287
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 53
$ws.Range("F4").Value = "Vedtrappmossa"
$ws.Range("G4").Value = "Crossocalyx hellerianus"
$ws.Range("H4").Value = "(Nees ex Lindenb.) Meyl."
$ws.Range("Q4").Value = 494363
$ws.Range("R4").Value = 6928873
$ws.Range("S4").Value = 30
# Row 5
$ws.Range("A5").Value = 111870913
$ws.Range("B5").Value = 90800
$ws.Range("D5").Value = "LC"
$ws.Range("E5").Value = 4364
$ws.Range("F5").Value = "Dropptaggsvamp"
$ws.Range("G5").Value = "Hydnellum ferrugineum"
$ws.Range("H5").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q5").Value = 494330
$ws.Range("R5").Value = 6928848
# Row 6
$ws.Range("A6").Value = 111870498
$ws.Range("B6").Value = 77636
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 6425
$ws.Range("F6").Value = "Garnlav"
$ws.Range("G6").Value = "Alectoria sarmentosa"
$ws.Range("H6").Value = "(Ach.) Ach."
$ws.Range("Q6").Value = 494301
$ws.Range("R6").Value = 6928922
$ws.Range("S6").Value = 20
# Row 7
$ws.Range("A7").Value = 111868975
$ws.Range("B7").Value = 90812
$ws.Range("D7").Value = "LC"
$ws.Range("E7").Value = 4366
$ws.Range("F7").Value = "Skarp dropptaggsvamp"
$ws.Range("G7").Value = "Hydnellum peckii"
$ws.Range("H7").Value = "Banker"
$ws.Range("Q7").Value = 494341
$ws.Range("R7").Value = 6928940
$ws.Range("S7").Value = 30
# Row 8
$ws.Range("A8").Value = 111869281
$ws.Range("B8").Value = 90804
$ws.Range("D8").Value = "VU"
$ws.Range("E8").Value = 4365
$ws.Range("F8").Value = "Smalfotad taggsvamp"
$ws.Range("G8").Value = "Hydnellum gracilipes"
$ws.Range("H8").Value = "(P.Karst) P.Karst"
$ws.Range("Q8").Value = 494333
$ws.Range("R8").Value = 6928943
$ws.Range("S8").Value = 30
$ws.Range("Z8").Value = "15:06"
$ws.Range("AB8").Value = "15:06"
# Row 9
$ws.Range("A9").Value = 111870723
$ws.Range("B9").Value = 77636
$ws.Range("Q9").Value = 494308
$ws.Range("R9").Value = 6928910
# Row 10
$ws.Range("A10").Value = 111870906
$ws.Range("B10").Value = 77636
$ws.Range("E10").Value = 6425
$ws.Range("F10").Value = "Garnlav"
$ws.Range("G10").Value = "Alectoria sarmentosa"
$ws.Range("H10").Value = "(Ach.) Ach."
$ws.Range("Q10").Value = 494330
$ws.Range("R10").Value = 6928848
$ws.Range("S10").Value = 20
$ws.Range("Z10").Value = "14:23"
$ws.Range("AB10").Value = "14:23"
# Row 11
$ws.Range("A11").Value = 111869523
$ws.Range("B11").Value = 56575
$ws.Range("D11").Value = "NT"
$ws.Range("E11").Value = 103021
$ws.Range("F11").Value = "Talltita"
$ws.Range("G11").Value = "Poecile montanus"
$ws.Range("H11").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("Q11").Value = 494333
$ws.Range("R11").Value = 6928943
$ws.Range("Z11").Value = "15:06"
$ws.Range("AB11").Value = "15:06"
# Row 12
$ws.Range("A12").Value = 111870057
$ws.Range("B12").Value = 90844
$ws.Range("E12").Value = 5449
$ws.Range("F12").Value = "Svart taggsvamp"
$ws.Range("G12").Value = "Phellodon niger"
$ws.Range("H12").Value = "(Fr.:Fr.) P.Karst."
$ws.Range("Q12").Value = 494314
$ws.Range("R12").Value = 6928937
# Row 13
$ws.Range("A13").Value = 111870880
$ws.Range("B13").Value = 90816
$ws.Range("D13").Value = "NT"
$ws.Range("E13").Value = 2059
$ws.Range("F13").Value = "Skrovlig taggsvamp"
$ws.Range("G13").Value = "Hydnellum scabrosum"
$ws.Range("H13").Value = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
# Row 14
$ws.Range("A14").Value = 111868438
$ws.Range("B14").Value = 89503
$ws.Range("D14").Value = "LC"
$ws.Range("E14").Value = 5447
$ws.Range("F14").Value = "Vedticka"
$ws.Range("G14").Value = "Fuscoporia viticola"
$ws.Range("H14").Value = "(Schwein.) Murrill"
$ws.Range("Q14").Value = 494363
$ws.Range("R14").Value = 6928873
$ws.Range("Z14").Value = "14:23"
$ws.Range("AB14").Value = "14:23"
# Row 15
$ws.Range("A15").Value = 111868497
$ws.Range("B15").Value = 90800
$ws.Range("P15").Value = "Motjärnen (Motjärnen), Jmt"
$ws.Range("Q15").Value = 494354
$ws.Range("R15").Value = 6928891
$ws.Range("S15").Value = 30
$ws.Range("J15").ClearContents()
$ws.Range("N15").ClearContents()
$ws.Range("AF15").ClearContents()
# Row 16
$ws.Range("A16").Value = 111868823
$ws.Range("B16").Value = 90800
$ws.Range("D16").Value = "LC"
$ws.Range("E16").Value = 4364
$ws.Range("F16").Value = "Dropptaggsvamp"
$ws.Range("G16").Value = "Hydnellum ferrugineum"
$ws.Range("H16").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("P16").Value = "Kläppberget, Kläppberget, Haverö, Jmt"
$ws.Range("Q16").Value = 494338
$ws.Range("R16").Value = 6928937
$ws.Range("S16").Value = 25
# Row 17
$ws.Range("B17").Value = 89902
